# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.152.95'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.837.37'
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '362.00'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +6.30%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '113.05'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.89%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.573'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +3.52%  '
$ws.Range('E8').Value = '  +0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.605'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +4.36%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '41.47'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.35%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0862'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  +1.21%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '20.01'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.33%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '7.81'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('D15').Value = '3.284.18'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = '2.841.01'
$ws.Range('E16').Value = '  +1.51%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.908'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +2.73%  '
$ws.Range('D18').Value = '52.045.99'
$ws.Range('E18').Value = '  -0.01%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.57'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +8.57%  '
$ws.Range('E20').Value = '  -1.42%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.53'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('D22').Value = '0.0₃0994'
$ws.Range('E22').Value = '  +1.32%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '70.37'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.32%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '268.11'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -3.58%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.84'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.64%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '27.10'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('E29').Value = '  +1.28%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0496'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +32.34%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '54.03'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +7.19%  '
$ws.Range('E32').Value = '  -1.93%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '35.14'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.57%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.86'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +2.15%  '
$ws.Range('E35').Value = '  +10.73%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0846'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.07'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.26'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.63%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '18.44'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -2.65%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.117'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '23.70'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.23%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '127.59'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.99%  '
$ws.Range('E44').Value = '  -7.16%  '
$ws.Range('E45').Value = '  -2.32%  '
$ws.Range('E46').Value = '  +2.87%  '
$ws.Range('D47').Value = '2.115.77'
$ws.Range('E47').Value = '  +1.01%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +12.86%  '
$ws.Range('E50').Value = '  +5.33%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '62.28'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +4.11%  '
